$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.467507333333333
$ws.Range("H2").Value = 4.402521999999999
$ws.Range("I2").Value = 0.1890754490804
$ws.Range("J2").Value = 0.1890754490804
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.618716333333334
$ws.Range("N2").Value = 7.856149000000001
$ws.Range("O2").Value = 0.07115908183301342
$ws.Range("P2").Value = 0.07115908183301341
$ws.Range("Q2").Value = 3.842985423086445
$ws.Range("R2").Value = 34.586868807778
$ws.Range("S2").Value = 0.01345443535372595
$ws.Range("T2").Value = 0.01345443535372594
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.467507333333333
$ws.Range("H3").Value = 4.402521999999999
$ws.Range("I3").Value = 0.1890754490804
$ws.Range("J3").Value = 0.1890754490804
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 15.503283
$ws.Range("N3").Value = 46.509849
$ws.Range("O3").Value = 0.4212748702999519
$ws.Range("P3").Value = 0.4212748702999519
$ws.Range("Q3").Value = 22.751181493242
$ws.Range("R3").Value = 204.760633439178
$ws.Range("S3").Value = 0.07965273528825068
$ws.Range("T3").Value = 0.07965273528825066
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.467507333333333
$ws.Range("H4").Value = 4.402521999999999
$ws.Range("I4").Value = 0.1890754490804
$ws.Range("J4").Value = 0.1890754490804
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.67887366666666
$ws.Range("N4").Value = 56.036621
$ws.Range("O4").Value = 0.5075660478670347
$ws.Range("P4").Value = 0.5075660478670347
$ws.Range("Q4").Value = 27.41138408424021
$ws.Range("R4").Value = 246.702456758162
$ws.Range("S4").Value = 0.0959682784384234
$ws.Range("T4").Value = 0.0959682784384234
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.293983333333333
$ws.Range("H5").Value = 18.88195
$ws.Range("I5").Value = 0.8109245509196
$ws.Range("J5").Value = 0.8109245509195999
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.618716333333334
$ws.Range("N5").Value = 7.856149000000001
$ws.Range("O5").Value = 0.07115908183301342
$ws.Range("P5").Value = 0.07115908183301341
$ws.Range("Q5").Value = 16.48215695672778
$ws.Range("R5").Value = 148.33941261055
$ws.Range("S5").Value = 0.05770464647928748
$ws.Range("T5").Value = 0.05770464647928746
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.293983333333333
$ws.Range("H6").Value = 18.88195
$ws.Range("I6").Value = 0.8109245509196
$ws.Range("J6").Value = 0.8109245509195999
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 15.503283
$ws.Range("N6").Value = 46.509849
$ws.Range("O6").Value = 0.4212748702999519
$ws.Range("P6").Value = 0.4212748702999519
$ws.Range("Q6").Value = 97.57740481395001
$ws.Range("R6").Value = 878.19664332555
$ws.Range("S6").Value = 0.3416221350117012
$ws.Range("T6").Value = 0.3416221350117012
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.293983333333333
$ws.Range("H7").Value = 18.88195
$ws.Range("I7").Value = 0.8109245509196
$ws.Range("J7").Value = 0.8109245509195999
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.67887366666666
$ws.Range("N7").Value = 56.036621
$ws.Range("O7").Value = 0.5075660478670347
$ws.Range("P7").Value = 0.5075660478670347
$ws.Range("Q7").Value = 117.5645195434389
$ws.Range("R7").Value = 1058.08067589095
$ws.Range("S7").Value = 0.4115977694286114
$ws.Range("T7").Value = 0.4115977694286113
